$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.018.92"
$ws.Range("E2").Value = "  +0.88%  "

$ws.Range("D3").Value = "'1.901.36"
$ws.Range("E3").Value = "  +0.65%  "

$ws.Range("D4").Value = "'0.9986"
$ws.Range("E4").Value = "  -0.23%  "

$ws.Range("D5").Value = "'0.7922"
$ws.Range("E5").Value = "  +0.03%  "

$ws.Range("D6").Value = "'244.64"
$ws.Range("E6").Value = "  +1.49%  "

$ws.Range("D7").Value = "'0.9991"
$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").Value = "'0.3165"
$ws.Range("E8").Value = "  +0.48%  "

$ws.Range("D9").Value = "'25.93"
$ws.Range("E9").Value = "  +1.71%  "

$ws.Range("D10").Value = "'0.07326"
$ws.Range("E10").Value = "  +4.83%  "

$ws.Range("D11").Value = "'0.08126"
$ws.Range("E11").Value = "  +1.20%  "

$ws.Range("D12").Value = "'0.7803"
$ws.Range("E12").Value = "  +3.14%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.984.78"
$ws.Range("E13").Value = "  +4.22%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.527"
$ws.Range("E14").Value = "  +4.59%  "

$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "'94.40"
$ws.Range("E15").Value = "  +2.46%  "

$ws.Range("D16").Value = "'6.271"
$ws.Range("E16").Value = "  +6.24%  "

$ws.Range("D17").Value = "'29.919.07"
$ws.Range("E17").Value = "  +0.45%  "

$ws.Range("D18").Value = "'14.05"
$ws.Range("E18").Value = "  +2.17%  "

$ws.Range("D19").Value = "'247.96"
$ws.Range("E19").Value = "  +1.78%  "

$ws.Range("D20").Value = "'0.000007856"
$ws.Range("E20").Value = "  +2.67%  "

$ws.Range("D21").Value = "'8.187"
$ws.Range("E21").Value = "  +1.14%  "

$ws.Range("E22").Value = "  -0.26%  "

$ws.Range("B23").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C23").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D23").Value = "'2.089.14"
$ws.Range("E23").Value = "  -3.38%  "

$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").Value = "'0.9986"
$ws.Range("E24").Value = "  -0.26%  "

$ws.Range("D25").Value = "'0.1614"
$ws.Range("E25").Value = "  -1.63%  "

$ws.Range("D26").Value = "'9.511"
$ws.Range("E26").Value = "  +2.68%  "

$ws.Range("D27").Value = "'163.48"
$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("D28").Value = "'18.87"
$ws.Range("E28").Value = "  +1.57%  "

$ws.Range("D29").Value = "'2.048"
$ws.Range("E29").Value = "  +0.34%  "

$ws.Range("D30").Value = "'1.445"
$ws.Range("E30").Value = "  +4.82%  "

$ws.Range("E31").Value = "  +1.09%  "

$ws.Range("D32").Value = "'4.508"
$ws.Range("E32").Value = "  +3.11%  "

$ws.Range("D33").Value = "'0.05637"
$ws.Range("E33").Value = "  -0.43%  "

$ws.Range("D34").Value = "'4.109"
$ws.Range("E34").Value = "  +1.59%  "

$ws.Range("D35").Value = "'1.256"
$ws.Range("E35").Value = "  -0.26%  "

$ws.Range("D36").Value = "'0.7578"
$ws.Range("E36").Value = "  +3.53%  "

$ws.Range("D37").Value = "'0.9984"
$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("D38").Value = "'2.672"
$ws.Range("E38").Value = "  +3.28%  "

$ws.Range("D39").Value = "'0.01941"
$ws.Range("E39").Value = "  +2.26%  "

$ws.Range("D40").Value = "'2.802"
$ws.Range("E40").Value = "  +0.96%  "

$ws.Range("D41").Value = "'1.146.42"
$ws.Range("E41").Value = "  +12.59%  "

$ws.Range("D42").Value = "'0.4490"
$ws.Range("E42").Value = "  +2.49%  "

$ws.Range("D43").Value = "'74.70"
$ws.Range("E43").Value = "  +3.55%  "

$ws.Range("D44").Value = "'5.988"
$ws.Range("E44").Value = "  +3.17%  "

$ws.Range("D45").Value = "'0.8589"
$ws.Range("E45").Value = "  +2.53%  "

$ws.Range("D46").Value = "'1.911"
$ws.Range("E46").Value = "  +3.33%  "

$ws.Range("D47").Value = "'3.179"
$ws.Range("E47").Value = "  +10.08%  "

$ws.Range("E48").Value = "  -0.17%  "

$ws.Range("D49").Value = "'102.27"
$ws.Range("E49").Value = "  +0.06%  "

$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "'7.570"
$ws.Range("E50").Value = "  +1.99%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'9.830"
$ws.Range("E51").Value = "  +0.00%  "
